$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column AQ: header "02-ago" plus values for rows 2-18
$ws.Range("AQ1").Value = "02-ago"

$ws.Range("AQ2").Value = 0
$ws.Range("AQ3").Value = 16.651895630789305
$ws.Range("AQ4").Value = 16.382439903877927
$ws.Range("AQ5").Value = 21.437824827707033
$ws.Range("AQ6").Value = 0
$ws.Range("AQ7").Value = 14.037490135098141
$ws.Range("AQ8").Value = 8.7223461005345975
$ws.Range("AQ9").Value = 14.279971016631979
$ws.Range("AQ10").Value = 15.161708685889922
$ws.Range("AQ11").Value = 12.435248895912705
$ws.Range("AQ12").Value = 0
$ws.Range("AQ13").Value = 10.016256522474805
$ws.Range("AQ14").Value = 0
$ws.Range("AQ15").Value = 0
$ws.Range("AQ16").Value = 12.697432316260828
$ws.Range("AQ17").Value = 0
$ws.Range("AQ18").Value = 0
